$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Deliver to Name" (M1) and "Remarks" (N1) are replaced by two new
# "Destination Code" / "Destination Name" columns, and "Remarks" moves
# to the new last column (O1).
$ws.Range("M1").Value = "Destination Code"
$ws.Range("N1").Value = "Destination Name"
$ws.Range("O1").Value = "Remarks"

# Match the header formatting (bold) used by the rest of row 1.
$ws.Range("O1").Font.Bold = $true

# Widen the affected / new columns to fit their new header text.
$ws.Columns.Item(13).ColumnWidth = 15.6667
$ws.Columns.Item(14).ColumnWidth = 16.5
$ws.Columns.Item(15).ColumnWidth = 7.6667

# Update the active selection to match the saved view.
$ws.Range("J5").Select() | Out-Null
